$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = 44495
$ws.Range("J5").Value = 200
$ws.Range("D6").Value = 44305
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 2500
$ws.Range("L6").Value = 2500
$ws.Range("M6").Value = 2500
$ws.Range("N6").Value = '$/unidad'
$ws.Range("P6").Value = 2500
$ws.Range("D7").Value = 44194
$ws.Range("I7").Value = 'Extra'
$ws.Range("J7").Value = 120
$ws.Range("K7").Value = 3500
$ws.Range("L7").Value = 3500
$ws.Range("M7").Value = 3500
$ws.Range("O7").Value = 'Región de O''Higgins'
$ws.Range("P7").Value = 3500
$ws.Range("D8").Value = 44194
$ws.Range("I8").Value = 'Primera'
$ws.Range("J8").Value = 200
$ws.Range("K8").Value = 3000
$ws.Range("L8").Value = 3000
$ws.Range("M8").Value = 3000
$ws.Range("P8").Value = 3000
$ws.Range("D9").Value = 44488
$ws.Range("J9").Value = 150
$ws.Range("K9").Value = 800
$ws.Range("L9").Value = 800
$ws.Range("M9").Value = 800
$ws.Range("N9").Value = '$/kilo (volumen en unidades)'
$ws.Range("O9").Value = 'Perú'
$ws.Range("P9").Value = 800
$ws.Range("D10").Value = 44483
$ws.Range("H10").Value = 'Sin especificar'
$ws.Range("I10").Value = 'Primera'
$ws.Range("J10").Value = 120
$ws.Range("K10").Value = 800
$ws.Range("L10").Value = 800
$ws.Range("M10").Value = 800
$ws.Range("N10").Value = '$/kilo (volumen en unidades)'
$ws.Range("O10").Value = 'Perú'
$ws.Range("P10").Value = 800
$ws.Range("I11").Value = 'Extra'
$ws.Range("J11").Value = 340
$ws.Range("K11").Value = 2500
$ws.Range("L11").Value = 2500
$ws.Range("M11").Value = 2500
$ws.Range("P11").Value = 2500
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 400
$ws.Range("K12").Value = 2000
$ws.Range("L12").Value = 2000
$ws.Range("M12").Value = 2000
$ws.Range("P12").Value = 2000
$ws.Range("I13").Value = 'Segunda'
$ws.Range("J13").Value = 300
$ws.Range("K13").Value = 1500
$ws.Range("L13").Value = 1500
$ws.Range("M13").Value = 1500
$ws.Range("P13").Value = 1500
$ws.Range("D14").Value = 44223
$ws.Range("H14").Value = 'Americana O Klondike'
$ws.Range("I14").Value = 'Tercera'
$ws.Range("J14").Value = 160
$ws.Range("K14").Value = 1000
$ws.Range("L14").Value = 1000
$ws.Range("M14").Value = 1000
$ws.Range("O14").Value = 'Región de O''Higgins'
$ws.Range("P14").Value = 1000
$ws.Range("D15").Value = 44477
$ws.Range("J15").Value = 80
$ws.Range("D16").Value = 44491
$ws.Range("J16").Value = 150
$ws.Range("D17").Value = 44497
$ws.Range("J17").Value = 250
$ws.Range("D18").Value = 44504
$ws.Range("J18").Value = 200
$ws.Range("D19").Value = 44510
$ws.Range("J19").Value = 250
$ws.Range("D20").Value = 44312
$ws.Range("J20").Value = 180
$ws.Range("K20").Value = 2500
$ws.Range("L20").Value = 2500
$ws.Range("M20").Value = 2500
$ws.Range("N20").Value = '$/unidad'
$ws.Range("P20").Value = 2500
$ws.Range("D21").Value = 44217
$ws.Range("J21").Value = 400
$ws.Range("K21").Value = 2500
$ws.Range("L21").Value = 2500
$ws.Range("M21").Value = 2500
$ws.Range("P21").Value = 2500
$ws.Range("D22").Value = 44217
$ws.Range("J22").Value = 280
$ws.Range("K22").Value = 2000
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = 2000
$ws.Range("P22").Value = 2000
